$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 439
$ws.Range("F6").Value = 1285
$ws.Range("F8").Value = 7590
$ws.Range("F9").Value = 91
$ws.Range("F12").Value = 8223
$ws.Range("F13").Value = 2
$ws.Range("F16").Value = 5627
$ws.Range("F17").Value = 5627
$ws.Range("F19").Value = 2579
$ws.Range("F20").Value = 1114
$ws.Range("F21").Value = 4586
$ws.Range("F22").Value = 334
$ws.Range("F25").Value = 29
$ws.Range("F26").Value = 503
$ws.Range("F27").Value = 3207
$ws.Range("F28").Value = 3207
$ws.Range("F30").Value = 11
$ws.Range("F31").Value = 2873
$ws.Range("F32").Value = 2873
$ws.Range("F33").Value = 23
$ws.Range("F34").Value = 329
$ws.Range("F35").Value = 123
$ws.Range("F36").Value = 285
$ws.Range("F37").Value = 3
$ws.Range("F38").Value = 640
$ws.Range("F39").Value = 11
$ws.Range("F41").Value = 1640
$ws.Range("F44").Value = 6
$ws.Range("F45").Value = 2657
$ws.Range("F47").Value = 2273
$ws.Range("F49").Value = 26

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 109
$ws.Range("F6").Value = 34
$ws.Range("F8").Value = 107

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 264
$ws.Range("F3").Value = 1314

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 264
$ws.Range("F4").Value = 1314
$ws.Range("F6").Value = 1285
$ws.Range("F7").Value = 7590
$ws.Range("F8").Value = 91
$ws.Range("F11").Value = 8223
$ws.Range("F14").Value = 5627
$ws.Range("F15").Value = 5627
$ws.Range("F17").Value = 2579
$ws.Range("F18").Value = 1115
$ws.Range("F19").Value = 4586
$ws.Range("F20").Value = 397
$ws.Range("F23").Value = 29
$ws.Range("F24").Value = 109
$ws.Range("F25").Value = 503
$ws.Range("F26").Value = 3207
$ws.Range("F27").Value = 3207
$ws.Range("F29").Value = 11
$ws.Range("F30").Value = 2873
$ws.Range("F31").Value = 2873
$ws.Range("F32").Value = 329
$ws.Range("F33").Value = 123
$ws.Range("F34").Value = 285
$ws.Range("F36").Value = 3
$ws.Range("F37").Value = 640
$ws.Range("F39").Value = 11
$ws.Range("F41").Value = 34
$ws.Range("F42").Value = 1640
$ws.Range("F45").Value = 6
$ws.Range("F46").Value = 2657
$ws.Range("F48").Value = 2273
$ws.Range("F50").Value = 26
$ws.Range("F52").Value = 107
